$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the metric/dispersion values in row 9 and row 10
$ws.Range("B9").Value = "r"
$ws.Range("J9").Value = "r"
$ws.Range("B10").Value = "g"
$ws.Range("J10").Value = "g"

# Update the active selection to O11
$ws.Range("O11").Select()
